# torquecurves.xlsx edit script
# - Rename "Racetech Racing Team TU" -> "Racetech"
# - Rename "Global Formula Racing(OREGON)" -> "GFR"
# - Insert a new blank worksheet ("Sheet1") before "Delft"
# - Add gearing data (G1-G6 / Primary / Final ratios) + a zero-speed
#   baseline row (row 3) to both the Racetech and GFR torque-curve sheets
# - Restore the various sheet-view selections / active tab left behind
#   by the author while doing the above

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Sheet renames
# ---------------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Racetech Racing Team TU")
$wsRace.Name = "Racetech"

$wsGFR = $wb.Worksheets.Item("Global Formula Racing(OREGON)")
$wsGFR.Name = "GFR"

# ---------------------------------------------------------------------
# 2. Insert a new empty sheet right before "Delft"
# ---------------------------------------------------------------------
$wsDelft = $wb.Worksheets.Item("Delft")
$wsNew = $wb.Worksheets.Add($wsDelft)

# ---------------------------------------------------------------------
# 3. Gearing data additions.
#    Shared-string order matters: write the text labels (G1..G6,
#    Primary, Final) for BOTH sheets first (GFR then Racetech) so the
#    new shared strings land in that exact order, and only write the
#    "-" placeholder text afterwards.
# ---------------------------------------------------------------------
$wsGFR.Range("A4").Value = "G1"
$wsGFR.Range("A5").Value = "G2"
$wsGFR.Range("A6").Value = "G3"
$wsGFR.Range("A7").Value = "G4"
$wsGFR.Range("A8").Value = "G5"
$wsGFR.Range("A9").Value = "G6"
$wsGFR.Range("A10").Value = "Primary"
$wsGFR.Range("A11").Value = "Final"

$wsRace.Range("A4").Value = "G1"
$wsRace.Range("A5").Value = "G2"
$wsRace.Range("A6").Value = "G3"
$wsRace.Range("A7").Value = "G4"
$wsRace.Range("A8").Value = "G5"
$wsRace.Range("A9").Value = "G6"
$wsRace.Range("A10").Value = "Primary"
$wsRace.Range("A11").Value = "Final"

# Numeric gear ratios - GFR
$wsGFR.Range("B4").Value = 2.846
$wsGFR.Range("B5").Value = 1.947
$wsGFR.Range("B6").Value = 1.556
$wsGFR.Range("B7").Value = 1.333
$wsGFR.Range("B8").Value = 1.19
$wsGFR.Range("B9").Value = 1.083
$wsGFR.Range("B10").Value = 1.955
$wsGFR.Range("B11").Value = 3

# Numeric gear ratios - Racetech (B9 is a text placeholder "-")
$wsRace.Range("B4").Value = 2.75
$wsRace.Range("B5").Value = 2
$wsRace.Range("B6").Value = 1.666
$wsRace.Range("B7").Value = 1.444
$wsRace.Range("B8").Value = 1.304
$wsRace.Range("B9").Value = "-"
$wsRace.Range("B10").Value = 2.111
$wsRace.Range("B11").Value = 2.625

# Zero-speed baseline row (row 3) on both torque-curve sheets
$wsGFR.Range("C3").Value = 0
$wsGFR.Range("D3").Value = 0
$wsGFR.Range("E3").Value = 2
$wsGFR.Range("F3").Formula = "=(C3*E3)/5252"

$wsRace.Range("C3").Value = 0
$wsRace.Range("D3").Value = 0
$wsRace.Range("E3").Value = 2
$wsRace.Range("F3").Formula = "=(C3*E3)/5252"

# ---------------------------------------------------------------------
# 4. Restore sheet-view / selection state left by the editing session
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Rennteam Uni Stuttgart")
$ws1.Activate()
$ws1.Range("H45").Select()

$ws3 = $wb.Worksheets.Item("Unicorn Race Engineering")
$ws3.Activate()
$ws3.Range("B11").Select()

$ws12 = $wb.Worksheets.Item("BA Motors(Berlin)")
$ws12.Activate()
$ws12.Range("B44").Select()

$ws8 = $wb.Worksheets.Item("WHZ Racing Team")
$ws8.Activate()
$ws8.Range("A29").Select()

$wsGFR.Activate()
$wsGFR.Range("G15").Select()

# Racetech is activated last so it ends up as the active/selected tab
$wsRace.Activate()
$wsRace.Range("B15").Select()

Write-Host "edit complete"
